$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 71430390
$ws.Range("I86").Value = 80248780
$ws.Range("K86").Value = 80248780
$ws.Range("M86").Value = -80247657
$ws.Range("H89").Value = 71430390
$ws.Range("I89").Value = 80248780
$ws.Range("K89").Value = 401243900
$ws.Range("M89").Value = -401238284
$ws.Range("H100").Value = 2218.8572
$ws.Range("I100").Value = 1512.375
$ws.Range("K100").Value = 1512.375
$ws.Range("M100").Value = -971.375
$ws.Range("H106").Value = 111112920
$ws.Range("I106").Value = 111112920
$ws.Range("K106").Value = 111112920
$ws.Range("M106").Value = -111112289
$ws.Range("H111").Value = 15631994
$ws.Range("I111").Value = 17864140
$ws.Range("K111").Value = 53592420
$ws.Range("M111").Value = -53589353
$ws.Range("H113").Value = 33340888
$ws.Range("I113").Value = 3129
$ws.Range("J113").Value = 62511428
$ws.Range("K113").Value = 3129
$ws.Range("L113").Value = 62511428
$ws.Range("M113").Value = 125
$ws.Range("N113").Value = -62517936
$ws.Range("H116").Value = 16671127
$ws.Range("I116").Value = 41668084
$ws.Range("J116").Value = 6490.5557
$ws.Range("K116").Value = 41668084
$ws.Range("L116").Value = 6490.5557
$ws.Range("M116").Value = -41664642
$ws.Range("N116").Value = -13374.5557
$ws.Range("H118").Value = 3500.889
$ws.Range("I118").Value = 3827.25
$ws.Range("J118").Value = 890
$ws.Range("K118").Value = 11481.75
$ws.Range("L118").Value = 2670
$ws.Range("M118").Value = -9824.75
$ws.Range("N118").Value = -5984
$ws.Range("H121").Value = 20000
$ws.Range("J121").Value = 20000
$ws.Range("L121").Value = 60000
$ws.Range("N121").Value = -63494
$ws.Range("H138").Value = 4714.298
$ws.Range("I138").Value = 881.7143
$ws.Range("K138").Value = 2645.1429
$ws.Range("M138").Value = 2494.8571

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2502674.2
$ws.Range("I32").Value = 2719883.5
$ws.Range("K32").Value = 2719883.5
$ws.Range("M32").Value = -2719596.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 43482116
$ws.Range("I86").Value = 2927.5557
$ws.Range("K86").Value = 2927.5557
$ws.Range("M86").Value = -1804.5557
$ws.Range("H89").Value = 43482116
$ws.Range("I89").Value = 2927.5557
$ws.Range("K89").Value = 14637.7785
$ws.Range("M89").Value = -9021.7785
$ws.Range("H107").Value = 41669100
$ws.Range("I107").Value = 48914900
$ws.Range("K107").Value = 48914900
$ws.Range("M107").Value = -48912980
$ws.Range("H134").Value = 3925.3416
$ws.Range("I134").Value = 2355.6323
$ws.Range("J134").Value = 11549.643
$ws.Range("K134").Value = 7066.896900000001
$ws.Range("L134").Value = 34648.929
$ws.Range("M134").Value = -4531.896900000001
$ws.Range("N134").Value = -39718.929

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7361.361
$ws.Range("I31").Value = 3275.5334
$ws.Range("K31").Value = 3275.5334
$ws.Range("M31").Value = -2980.5334
$ws.Range("H34").Value = 7361.361
$ws.Range("I34").Value = 3275.5334
$ws.Range("K34").Value = 3275.5334
$ws.Range("M34").Value = -3073.5334
$ws.Range("H87").Value = 67330
$ws.Range("J87").Value = 67330
$ws.Range("L87").Value = 67330
$ws.Range("N87").Value = -69702
$ws.Range("H90").Value = 67330
$ws.Range("J90").Value = 67330
$ws.Range("L90").Value = 201990
$ws.Range("N90").Value = -213846
$ws.Range("H99").Value = 7308.826
$ws.Range("I99").Value = 8679.5
$ws.Range("K99").Value = 8679.5
$ws.Range("M99").Value = -7181.5
$ws.Range("H126").Value = 7308.826
$ws.Range("I126").Value = 8679.5
$ws.Range("K126").Value = 26038.5
$ws.Range("M126").Value = -23568.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3638607.8
$ws.Range("I5").Value = 5714498
$ws.Range("J5").Value = 5799.75
$ws.Range("K5").Value = 17143494
$ws.Range("L5").Value = 17399.25
$ws.Range("M5").Value = -17143382
$ws.Range("N5").Value = -17623.25
$ws.Range("H135").Value = 3638607.8
$ws.Range("I135").Value = 5714498
$ws.Range("J135").Value = 5799.75
$ws.Range("K135").Value = 51430482
$ws.Range("L135").Value = 52197.75
$ws.Range("M135").Value = -51427947
$ws.Range("N135").Value = -57267.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6081.5654
$ws.Range("I70").Value = 4882.0645
$ws.Range("J70").Value = 8560.532999999999
$ws.Range("K70").Value = 4882.0645
$ws.Range("L70").Value = 8560.532999999999
$ws.Range("M70").Value = -4612.0645
$ws.Range("N70").Value = -9100.532999999999
$ws.Range("H73").Value = 6081.5654
$ws.Range("I73").Value = 4882.0645
$ws.Range("J73").Value = 8560.532999999999
$ws.Range("K73").Value = 4882.0645
$ws.Range("L73").Value = 8560.532999999999
$ws.Range("M73").Value = -3946.0645
$ws.Range("N73").Value = -10432.533
$ws.Range("H113").Value = 5796.9707
$ws.Range("I113").Value = 3327.4443
$ws.Range("J113").Value = 6686
$ws.Range("K113").Value = 3327.4443
$ws.Range("L113").Value = 6686
$ws.Range("M113").Value = -1157.4443
$ws.Range("N113").Value = -11026
$ws.Range("H122").Value = 44769.32
$ws.Range("I122").Value = 93808.63
$ws.Range("K122").Value = 281425.89
$ws.Range("M122").Value = -278975.89
$ws.Range("H126").Value = 2839.3572
$ws.Range("I126").Value = 2689.75
$ws.Range("J126").Value = 2899.2
$ws.Range("K126").Value = 8069.25
$ws.Range("L126").Value = 8697.599999999999
$ws.Range("M126").Value = -5599.25
$ws.Range("N126").Value = -13637.6
$ws.Range("H132").Value = 2544.3784
$ws.Range("I132").Value = 2507.3872
$ws.Range("K132").Value = 7522.1616
$ws.Range("M132").Value = -4992.1616

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1526.9166
$ws.Range("I16").Value = 1545.3182
$ws.Range("K16").Value = 1545.3182
$ws.Range("M16").Value = -1375.3182
$ws.Range("H46").Value = 12347289
$ws.Range("I46").Value = 896
$ws.Range("K46").Value = 896
$ws.Range("M46").Value = -708
$ws.Range("H82").Value = 2128.5715
$ws.Range("I82").Value = 1982.75
$ws.Range("K82").Value = 1982.75
$ws.Range("M82").Value = -1621.75
$ws.Range("H85").Value = 2128.5715
$ws.Range("I85").Value = 1982.75
$ws.Range("K85").Value = 1982.75
$ws.Range("M85").Value = -734.75
$ws.Range("H132").Value = 16674716
$ws.Range("I132").Value = 29416362
$ws.Range("K132").Value = 88249086
$ws.Range("M132").Value = -88246556

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1811.625
$ws.Range("I107").Value = 2664.6667
$ws.Range("J107").Value = 1299.8
$ws.Range("K107").Value = 7994.000100000001
$ws.Range("L107").Value = 3899.4
$ws.Range("M107").Value = -6074.000100000001
$ws.Range("N107").Value = -7739.4
$ws.Range("H113").Value = 2547.889
$ws.Range("I113").Value = 2371.1428
$ws.Range("K113").Value = 7113.428400000001
$ws.Range("M113").Value = -4943.428400000001
$ws.Range("H132").Value = 17264164
$ws.Range("I132").Value = 20839732
$ws.Range("K132").Value = 62519196
$ws.Range("M132").Value = -62516666
